# Add a "Sample Data" column (column I) with ten sample rows to the
# "May 2024" sheet, matching the "April 2024" sheet that already has it,
# and refresh the "April 2024" header cell / selection state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "April 2024" sheet already has the Sample Data column (I1:I11). Its
# header cell currently points at a style that is about to become
# unused/removed, so re-point it at the same bold-header style used by
# the rest of the header row (same style as A1:H1).
# ---------------------------------------------------------------------
$wsApril = $wb.Worksheets.Item("April 2024")
$wsApril.Activate()

$wsApril.Range("A1").Copy()
$wsApril.Range("I1").PasteSpecial(-4122)
$wsApril.Range("I1").Value = "Sample Data"

$wsApril.Range("B38:B39").Select()

# ---------------------------------------------------------------------
# "May 2024" sheet: add the same "Sample Data" column with ten rows of
# sample values, matching column I of "April 2024".
# ---------------------------------------------------------------------
$wsMay = $wb.Worksheets.Item("May 2024")
$wsMay.Activate()

$wsMay.Range("A1").Copy()
$wsMay.Range("I1").PasteSpecial(-4122)
$wsMay.Range("I1").Value = "Sample Data"

$wsMay.Range("I2").Value = "Data 1"
$wsMay.Range("I3").Value = "Data 2"
$wsMay.Range("I4").Value = "Data 3"
$wsMay.Range("I5").Value = "Data 4"
$wsMay.Range("I6").Value = "Data 5"
$wsMay.Range("I7").Value = "Data 6"
$wsMay.Range("I8").Value = "Data 7"
$wsMay.Range("I9").Value = "Data 8"
$wsMay.Range("I10").Value = "Data 9"
$wsMay.Range("I11").Value = "Data 10"

# Match the column width ("bestFit" width of 12 used on "April 2024"'s
# identical column).
$wsMay.Columns.Item(9).ColumnWidth = 11.1666666666667

# Leave "May 2024" the active sheet/tab with A1 selected.
$wsMay.Range("A1").Select()
